# Auto-generated edit script: updates leve-crafting profit figures
# across all 8 class sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# to reflect refreshed market-board prices from the scheduled runner.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1844.6
$ws.Range("J17").Value = 1844.6
$ws.Range("L17").Value = 5533.799999999999
$ws.Range("N17").Value = -5869.799999999999
$ws.Range("H28").Value = 1118.8077
$ws.Range("I28").Value = 1174.5834
$ws.Range("J28").Value = 449.5
$ws.Range("K28").Value = 1174.5834
$ws.Range("L28").Value = 449.5
$ws.Range("M28").Value = -689.5834
$ws.Range("N28").Value = -1419.5
$ws.Range("H34").Value = 8371.666999999999
$ws.Range("I34").Value = 8371.666999999999
$ws.Range("K34").Value = 8371.666999999999
$ws.Range("M34").Value = -8168.666999999999
$ws.Range("H36").Value = 8371.666999999999
$ws.Range("I36").Value = 8371.666999999999
$ws.Range("K36").Value = 8371.666999999999
$ws.Range("M36").Value = -7656.666999999999
$ws.Range("H51").Value = 6722.4287
$ws.Range("J51").Value = 7838.8335
$ws.Range("L51").Value = 7838.8335
$ws.Range("N51").Value = -8806.833500000001
$ws.Range("H137").Value = 337785.38
$ws.Range("I137").Value = 480860.28
$ws.Range("J137").Value = 3943.889
$ws.Range("K137").Value = 1442580.84
$ws.Range("L137").Value = 11831.667
$ws.Range("M137").Value = -1440030.84
$ws.Range("N137").Value = -16931.667
$ws.Range("H138").Value = 4307.023
$ws.Range("I138").Value = 1022
$ws.Range("K138").Value = 3066
$ws.Range("M138").Value = 2074

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13726.288
$ws.Range("I32").Value = 14013.203
$ws.Range("K32").Value = 14013.203
$ws.Range("M32").Value = -13726.203
$ws.Range("H43").Value = 23190
$ws.Range("J43").Value = 23190
$ws.Range("L43").Value = 23190
$ws.Range("N43").Value = -23816
$ws.Range("H61").Value = 4050.9285
$ws.Range("I61").Value = 2575.0527
$ws.Range("K61").Value = 2575.0527
$ws.Range("M61").Value = -2363.0527
$ws.Range("H74").Value = 1747.2
$ws.Range("I74").Value = 1181.2142
$ws.Range("K74").Value = 1181.2142
$ws.Range("M74").Value = -307.2141999999999
$ws.Range("H77").Value = 1747.2
$ws.Range("I77").Value = 1181.2142
$ws.Range("K77").Value = 5906.071
$ws.Range("M77").Value = -1538.071
$ws.Range("H132").Value = 14489.631
$ws.Range("I132").Value = 16389.81
$ws.Range("J132").Value = 6677.778
$ws.Range("K132").Value = 49169.43000000001
$ws.Range("L132").Value = 20033.334
$ws.Range("M132").Value = -46639.43000000001
$ws.Range("N132").Value = -25093.334
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 4050.9285
$ws.Range("I136").Value = 2575.0527
$ws.Range("K136").Value = 7725.158100000001
$ws.Range("M136").Value = -5175.158100000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4004.5
$ws.Range("I94").Value = 4009
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 4009
$ws.Range("L94").Value = 4000
$ws.Range("M94").Value = -3558
$ws.Range("N94").Value = -4902

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5204.9395
$ws.Range("I31").Value = 4165.923
$ws.Range("J31").Value = 5880.3
$ws.Range("K31").Value = 4165.923
$ws.Range("L31").Value = 5880.3
$ws.Range("M31").Value = -3870.923
$ws.Range("N31").Value = -6470.3
$ws.Range("H34").Value = 5204.9395
$ws.Range("I34").Value = 4165.923
$ws.Range("J34").Value = 5880.3
$ws.Range("K34").Value = 4165.923
$ws.Range("L34").Value = 5880.3
$ws.Range("M34").Value = -3963.923
$ws.Range("N34").Value = -6284.3
$ws.Range("H58").Value = 557154.6
$ws.Range("I58").Value = 1329.75
$ws.Range("K58").Value = 1329.75
$ws.Range("M58").Value = -1126.75
$ws.Range("H99").Value = 10893.706
$ws.Range("I99").Value = 10605.211
$ws.Range("J99").Value = 11259.134
$ws.Range("K99").Value = 10605.211
$ws.Range("L99").Value = 11259.134
$ws.Range("M99").Value = -9107.210999999999
$ws.Range("N99").Value = -14255.134
$ws.Range("H126").Value = 10893.706
$ws.Range("I126").Value = 10605.211
$ws.Range("J126").Value = 11259.134
$ws.Range("K126").Value = 31815.633
$ws.Range("L126").Value = 33777.402
$ws.Range("M126").Value = -29345.633
$ws.Range("N126").Value = -38717.402
$ws.Range("H132").Value = 7102762
$ws.Range("I132").Value = 7943180.5
$ws.Range("J132").Value = 43247.2
$ws.Range("K132").Value = 23829541.5
$ws.Range("L132").Value = 129741.6
$ws.Range("M132").Value = -23827011.5
$ws.Range("N132").Value = -134801.6
$ws.Range("H134").Value = 2529.8125
$ws.Range("J134").Value = 2416.7778
$ws.Range("L134").Value = 7250.3334
$ws.Range("N134").Value = -12320.3334
$ws.Range("H136").Value = 557154.6
$ws.Range("I136").Value = 1329.75
$ws.Range("K136").Value = 3989.25
$ws.Range("M136").Value = -1439.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 182.72414
$ws.Range("J2").Value = 234.07692
$ws.Range("L2").Value = 1404.46152
$ws.Range("N2").Value = -1630.46152
$ws.Range("H4").Value = 4002717.5
$ws.Range("I4").Value = 1793594.9
$ws.Range("K4").Value = 5380784.699999999
$ws.Range("M4").Value = -5380672.699999999
$ws.Range("H5").Value = 1170.5
$ws.Range("I5").Value = 776.5
$ws.Range("J5").Value = 1663
$ws.Range("K5").Value = 2329.5
$ws.Range("L5").Value = 4989
$ws.Range("M5").Value = -2217.5
$ws.Range("N5").Value = -5213
$ws.Range("H56").Value = 45462104
$ws.Range("I56").Value = 45462104
$ws.Range("K56").Value = 45462104
$ws.Range("M56").Value = -45461574
$ws.Range("H103").Value = 545.38464
$ws.Range("J103").Value = 686.375
$ws.Range("L103").Value = 2059.125
$ws.Range("N103").Value = -3817.125
$ws.Range("H113").Value = 838.4
$ws.Range("J113").Value = 843.25
$ws.Range("L113").Value = 2529.75
$ws.Range("N113").Value = -6869.75
$ws.Range("H131").Value = 2000457.2
$ws.Range("I131").Value = 112039.555
$ws.Range("J131").Value = 2459802
$ws.Range("K131").Value = 336118.665
$ws.Range("L131").Value = 7379406
$ws.Range("M131").Value = -331078.665
$ws.Range("N131").Value = -7389486
$ws.Range("H135").Value = 1170.5
$ws.Range("I135").Value = 776.5
$ws.Range("J135").Value = 1663
$ws.Range("K135").Value = 6988.5
$ws.Range("L135").Value = 14967
$ws.Range("M135").Value = -4453.5
$ws.Range("N135").Value = -20037
$ws.Range("H140").Value = 8620.620999999999
$ws.Range("I140").Value = 2381
$ws.Range("J140").Value = 24999.625
$ws.Range("K140").Value = 7143
$ws.Range("L140").Value = 74998.875
$ws.Range("M140").Value = -1963
$ws.Range("N140").Value = -85358.875

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 34500
$ws.Range("J95").Value = 34500
$ws.Range("L95").Value = 34500
$ws.Range("N95").Value = -39992
$ws.Range("H102").Value = 4670.5557
$ws.Range("I102").Value = 4123.0713
$ws.Range("K102").Value = 4123.0713
$ws.Range("M102").Value = -2501.0713
$ws.Range("H123").Value = 44925
$ws.Range("J123").Value = 44925
$ws.Range("L123").Value = 44925
$ws.Range("N123").Value = -49825
$ws.Range("H126").Value = 4297.92
$ws.Range("I126").Value = 2550
$ws.Range("K126").Value = 7650
$ws.Range("M126").Value = -5180
$ws.Range("H132").Value = 605209.25
$ws.Range("I132").Value = 122599.234
$ws.Range("K132").Value = 367797.702
$ws.Range("M132").Value = -365267.702

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2937.0476
$ws.Range("I46").Value = 805.6
$ws.Range("J46").Value = 3603.125
$ws.Range("K46").Value = 805.6
$ws.Range("L46").Value = 3603.125
$ws.Range("M46").Value = -617.6
$ws.Range("N46").Value = -3979.125
$ws.Range("H68").Value = 2930.56
$ws.Range("I68").Value = 2327.45
$ws.Range("K68").Value = 2327.45
$ws.Range("M68").Value = -1578.45
$ws.Range("H71").Value = 2930.56
$ws.Range("I71").Value = 2327.45
$ws.Range("K71").Value = 11637.25
$ws.Range("M71").Value = -7893.25
$ws.Range("H132").Value = 1853.17
$ws.Range("I132").Value = 1848.5851
$ws.Range("J132").Value = 1925
$ws.Range("K132").Value = 5545.7553
$ws.Range("L132").Value = 5775
$ws.Range("M132").Value = -3015.7553
$ws.Range("N132").Value = -10835
$ws.Range("H139").Value = 78707.5
$ws.Range("J139").Value = 78707.5
$ws.Range("L139").Value = 78707.5
$ws.Range("N139").Value = -88987.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3255.625
$ws.Range("I81").Value = 864.2857
$ws.Range("J81").Value = 19995
$ws.Range("K81").Value = 1728.5714
$ws.Range("L81").Value = 39990
$ws.Range("M81").Value = -667.5714
$ws.Range("N81").Value = -42112
$ws.Range("H84").Value = 3255.625
$ws.Range("I84").Value = 864.2857
$ws.Range("J84").Value = 19995
$ws.Range("K84").Value = 8642.857
$ws.Range("L84").Value = 199950
$ws.Range("M84").Value = -3338.857
$ws.Range("N84").Value = -210558
$ws.Range("H100").Value = 3432.9644
$ws.Range("I100").Value = 3479.65
$ws.Range("J100").Value = 3316.25
$ws.Range("K100").Value = 6959.3
$ws.Range("L100").Value = 6632.5
$ws.Range("M100").Value = -6418.3
$ws.Range("N100").Value = -7714.5
$ws.Range("H132").Value = 3683.1516
$ws.Range("I132").Value = 1150.32
$ws.Range("J132").Value = 11598.25
$ws.Range("K132").Value = 3450.96
$ws.Range("L132").Value = 34794.75
$ws.Range("M132").Value = -920.96
$ws.Range("N132").Value = -39854.75
$ws.Range("H133").Value = 60818.4
$ws.Range("J133").Value = 60818.4
$ws.Range("L133").Value = 60818.4
$ws.Range("N133").Value = -70938.39999999999
$ws.Range("H141").Value = 74163.57000000001
$ws.Range("J141").Value = 74163.57000000001
$ws.Range("L141").Value = 74163.57000000001
$ws.Range("N141").Value = -84523.57000000001
